# Auto-update draw results: append the 2025-12-11 Pick 3 draw as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Range("A86:E86")

# All columns in this sheet store plain text (dates/phase codes/results are
# kept as literal strings, not numbers/dates) - force text formatting before
# assigning so Excel doesn't auto-coerce "2025-12-11"/"251211" into a date
# serial / number, then clear the formatting again so the new cells end up
# with the same (default) style as every other data row.
$newRow.NumberFormat = "@"

$ws.Range("A86").Value = "2025-12-11"
$ws.Range("B86").Value = "Pick 3"
$ws.Range("C86").Value = "251211"
$ws.Range("D86").Value = "0-8-1"
$ws.Range("E86").Value = "2025-12-11T21:52:25.254+04:00"

$newRow.ClearFormats()
